{"js": "const replacements = [\n  [\"2024-02-15 Thursday\", \"2024-02-16 Friday\"],\n  [\"48\u00d792=4416\", \"19\u00d752=988\"],\n  [\"96\u00d741=3936\", \"27\u00d727=729\"],\n  [\"50\u00d780=4000\", \"17\u00d756=952\"],\n  [\"93\u00d730=2790\", \"44\u00d760=2640\"],\n  [\"34\u00d714=476\", \"63\u00d728=1764\"],\n  [\"57\u00d745=2565\", \"91\u00d754=4914\"],\n  [\"22\u00d772=1584\", \"85\u00d732=2720\"],\n  [\"19\u00d766=1254\", \"65\u00d765=4225\"],\n  [\"70\u00d742=2940\", \"64\u00d720=1280\"],\n  [\"73\u00d726=1898\", \"32\u00d718=576\"],\n  [\"92\u00d734=3128\", \"56\u00d743=2408\"],\n  [\"69\u00d722=1518\", \"69\u00d712=828\"],\n  [\"32\u00d764=2048\", \"39\u00d731=1209\"],\n  [\"89\u00d731=2759\", \"13\u00d781=1053\"],\n  [\"58\u00d793=5394\", \"39\u00d750=1950\"],\n  [\"52\u00d792=4784\", \"81\u00d758=4698\"],\n  [\"50\u00d733=1650\", \"36\u00d738=1368\"],\n  [\"95\u00d718=1710\", \"15\u00d760=900\"],\n  [\"54\u00d763=3402\", \"39\u00d722=858\"],\n  [\"18\u00d756=1008\", \"37\u00d757=2109\"],\n  [\"37\u00d743=1591\", \"98\u00d782=8036\"],\n  [\"99\u00d724=2376\", \"65\u00d773=4745\"],\n  [\"28\u00d797=2716\", \"12\u00d781=972\"],\n  [\"18\u00d784=1512\", \"13\u00d726=338\"],\n  [\"50\u00d783=4150\", \"82\u00d754=4428\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('Text not found: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-15 Thursday\", \"2024-02-16 Friday\"),\n    @(\"48\u00d792=4416\", \"19\u00d752=988\"),\n    @(\"96\u00d741=3936\", \"27\u00d727=729\"),\n    @(\"50\u00d780=4000\", \"17\u00d756=952\"),\n    @(\"93\u00d730=2790\", \"44\u00d760=2640\"),\n    @(\"34\u00d714=476\", \"63\u00d728=1764\"),\n    @(\"57\u00d745=2565\", \"91\u00d754=4914\"),\n    @(\"22\u00d772=1584\", \"85\u00d732=2720\"),\n    @(\"19\u00d766=1254\", \"65\u00d765=4225\"),\n    @(\"70\u00d742=2940\", \"64\u00d720=1280\"),\n    @(\"73\u00d726=1898\", \"32\u00d718=576\"),\n    @(\"92\u00d734=3128\", \"56\u00d743=2408\"),\n    @(\"69\u00d722=1518\", \"69\u00d712=828\"),\n    @(\"32\u00d764=2048\", \"39\u00d731=1209\"),\n    @(\"89\u00d731=2759\", \"13\u00d781=1053\"),\n    @(\"58\u00d793=5394\", \"39\u00d750=1950\"),\n    @(\"52\u00d792=4784\", \"81\u00d758=4698\"),\n    @(\"50\u00d733=1650\", \"36\u00d738=1368\"),\n    @(\"95\u00d718=1710\", \"15\u00d760=900\"),\n    @(\"54\u00d763=3402\", \"39\u00d722=858\"),\n    @(\"18\u00d756=1008\", \"37\u00d757=2109\"),\n    @(\"37\u00d743=1591\", \"98\u00d782=8036\"),\n    @(\"99\u00d724=2376\", \"65\u00d773=4745\"),\n    @(\"28\u00d797=2716\", \"12\u00d781=972\"),\n    @(\"18\u00d784=1512\", \"13\u00d726=338\"),\n    @(\"50\u00d783=4150\", \"82\u00d754=4428\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Host \"Not found: $oldText\"\n    }\n}"}
